$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> Clcf1/Cntfr -> ECs)
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1594223333333333
$ws.Range("H2").Value = 0.478267
$ws.Range("I2").Value = 0.01552338951653915
$ws.Range("J2").Value = 0.01552338951653915
$ws.Range("M2").Value = 0.01031333333333333
$ws.Range("N2").Value = 0.03094
$ws.Range("O2").Value = 0.001146416507271297
$ws.Range("P2").Value = 0.001146416507271297
$ws.Range("Q2").Value = 0.001644175664444444
$ws.Range("R2").Value = 0.01479758098
$ws.Range("S2").Value = 0.00001779626999056267
$ws.Range("T2").Value = 0.00001779626999056267

# Row 3 (ECs -> Clcf1/Cntfr -> FAPs)
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1594223333333333
$ws.Range("H3").Value = 0.478267
$ws.Range("I3").Value = 0.01552338951653915
$ws.Range("J3").Value = 0.01552338951653915
$ws.Range("O3").Value = 0.8303652607489888
$ws.Range("P3").Value = 0.8303652607489886
$ws.Range("Q3").Value = 1.190899071728444
$ws.Range("R3").Value = 10.718091645556
$ws.Range("S3").Value = 0.01289008338360915
$ws.Range("T3").Value = 0.01289008338360915

# Row 4 (ECs -> Clcf1/Cntfr -> MuSCs)
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1594223333333333
$ws.Range("H4").Value = 0.478267
$ws.Range("I4").Value = 0.01552338951653915
$ws.Range("J4").Value = 0.01552338951653915
$ws.Range("M4").Value = 1.515746
$ws.Range("N4").Value = 4.547238
$ws.Range("O4").Value = 0.16848832274374
$ws.Range("P4").Value = 0.16848832274374
$ws.Range("Q4").Value = 0.2416437640606667
$ws.Range("R4").Value = 2.174793876546
$ws.Range("S4").Value = 0.002615509862939439
$ws.Range("T4").Value = 0.002615509862939439

# Row 5 (FAPs -> Clcf1/Cntfr -> ECs)
$ws.Range("I5").Value = 0.1862883666449807
$ws.Range("J5").Value = 0.1862883666449807
$ws.Range("M5").Value = 0.01031333333333333
$ws.Range("N5").Value = 0.03094
$ws.Range("O5").Value = 0.001146416507271297
$ws.Range("P5").Value = 0.001146416507271297
$ws.Range("Q5").Value = 0.01973092272666667
$ws.Range("R5").Value = 0.17757830454
$ws.Range("S5").Value = 0.0002135640586344135
$ws.Range("T5").Value = 0.0002135640586344134

# Row 6 (FAPs -> Clcf1/Cntfr -> FAPs)
$ws.Range("I6").Value = 0.1862883666449807
$ws.Range("J6").Value = 0.1862883666449807
$ws.Range("O6").Value = 0.8303652607489888
$ws.Range("P6").Value = 0.8303652607489886
$ws.Range("S6").Value = 0.1546873881436626
$ws.Range("T6").Value = 0.1546873881436626

# Row 7 (FAPs -> Clcf1/Cntfr -> MuSCs)
$ws.Range("I7").Value = 0.1862883666449807
$ws.Range("J7").Value = 0.1862883666449807
$ws.Range("M7").Value = 1.515746
$ws.Range("N7").Value = 4.547238
$ws.Range("O7").Value = 0.16848832274374
$ws.Range("P7").Value = 0.16848832274374
$ws.Range("Q7").Value = 2.899844912662001
$ws.Range("R7").Value = 26.098604213958
$ws.Range("S7").Value = 0.03138741444268368
$ws.Range("T7").Value = 0.03138741444268368

# Row 8 (MuSCs -> Clcf1/Cntfr -> ECs)
$ws.Range("G8").Value = 8.197245333333333
$ws.Range("H8").Value = 24.591736
$ws.Range("I8").Value = 0.7981882438384801
$ws.Range("J8").Value = 0.7981882438384801
$ws.Range("M8").Value = 0.01031333333333333
$ws.Range("N8").Value = 0.03094
$ws.Range("O8").Value = 0.001146416507271297
$ws.Range("P8").Value = 0.001146416507271297
$ws.Range("Q8").Value = 0.08454092353777777
$ws.Range("R8").Value = 0.7608683118399999
$ws.Range("S8").Value = 0.0009150561786463204
$ws.Range("T8").Value = 0.0009150561786463203

# Row 9 (MuSCs -> Clcf1/Cntfr -> FAPs)
$ws.Range("G9").Value = 8.197245333333333
$ws.Range("H9").Value = 24.591736
$ws.Range("I9").Value = 0.7981882438384801
$ws.Range("J9").Value = 0.7981882438384801
$ws.Range("O9").Value = 0.8303652607489888
$ws.Range("P9").Value = 0.8303652607489886
$ws.Range("Q9").Value = 61.23415492724978
$ws.Range("R9").Value = 551.1073943452479
$ws.Range("S9").Value = 0.6627877892217169
$ws.Range("T9").Value = 0.6627877892217169

# Row 10 (MuSCs -> Clcf1/Cntfr -> MuSCs)
$ws.Range("G10").Value = 8.197245333333333
$ws.Range("H10").Value = 24.591736
$ws.Range("I10").Value = 0.7981882438384801
$ws.Range("J10").Value = 0.7981882438384801
$ws.Range("M10").Value = 1.515746
$ws.Range("N10").Value = 4.547238
$ws.Range("O10").Value = 0.16848832274374
$ws.Range("P10").Value = 0.16848832274374
$ws.Range("Q10").Value = 12.42494182501867
$ws.Range("R10").Value = 111.824476425168
$ws.Range("S10").Value = 0.1344853984381169
$ws.Range("T10").Value = 0.1344853984381169
